$wb = $excel.ActiveWorkbook

# --- Add the new "Antonio" worksheet after "Foglio1" ---
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "Antonio"

# --- Text labels (entered in this order so shared-string indices line up
#     the same way the original author's edit produced them) ---
$ws.Range("B3").Value = "camera"
$ws.Range("B2").Value = "installazione + materiale"
$ws.Range("D1").Value = "num"
$ws.Range("C1").Value = "singolo"
$ws.Range("B4").Value = "modem"
$ws.Range("B5").Value = "switch"
$ws.Range("B6").Value = "miniPC"
$ws.Range("B7").Value = "scheda semaforo??"

$euroFmt = "#,##0.00\ ""€"""

# --- Row 2: installazione + materiale ---
$ws.Range("C2").Value = 3000
$ws.Range("C2").NumberFormat = $euroFmt
$ws.Range("D2").Value = 1
$ws.Range("E2").Formula = "=D2*C2"
$ws.Range("E2").NumberFormat = $euroFmt

# --- Row 3: camera ---
$ws.Range("C3").Value = 320
$ws.Range("C3").NumberFormat = $euroFmt
$ws.Range("D3").Value = 4
$ws.Range("E3").Formula = "=D3*C3"
$ws.Range("E3").NumberFormat = $euroFmt

# --- Rows 4-6: modem / switch / miniPC, filled as one formula so the
#     engine emits a shared formula (t="shared") like the source file ---
$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 1
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = 1
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = 1
$ws.Range("E4:E6").Formula = "=D4*C4"

$ws.Range("C4").NumberFormat = $euroFmt
$ws.Range("C5").NumberFormat = $euroFmt
$ws.Range("C6").NumberFormat = $euroFmt
$ws.Range("E4").NumberFormat = $euroFmt
$ws.Range("E5").NumberFormat = $euroFmt
$ws.Range("E6").NumberFormat = $euroFmt

# --- Row 7: scheda semaforo?? (only the currency-formatted C cell, no values) ---
$ws.Range("C7").NumberFormat = $euroFmt

# --- Row 8: total ---
$ws.Range("E8").Formula = "=SUM(E2:E7)"
$ws.Range("E8").NumberFormat = $euroFmt

# --- Wrap text on the "installazione + materiale" label ---
$ws.Range("B2").WrapText = $true

# --- Column widths ---
$ws.Range("B2:D2").ColumnWidth = 12.140625
$ws.Range("E1").ColumnWidth = 9.85546875

# --- View: make Antonio the active tab / sheet, set zoom + selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("F13").Select()
